$d = $word.ActiveDocument

# Insert the new sentence about the dajaxice bug's error message right after
# the existing sentence "...details the bug, and the easy workaround."
# wildcard = $false, MatchCase = $true, WholeWord = $false, MatchSoundsLike = $false,
# MatchAllWordForms = $false, Forward = $true, Wrap = 1 (wdFindContinue),
# Format = $false, ReplaceWith = <text>, Replace = 2 (wdReplaceAll)
$d.Content.Find.Execute(
    "and the easy workaround.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "and the easy workaround. The error message when you try to load the site will tell you the path of the urls.py file that needs to be edited.",
    2
)
